$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REPORT")

# Duplicate the "Selection" review block (rows 32:43) down into rows 45:56,
# matching a copy/paste of that block further down the sheet.
# Columns A, B-D and F-H are copied as a straight block so the resulting
# row "spans" metadata matches a real Excel paste. Column E (the date column)
# is populated separately with Value2 + a formats-only paste so we reuse the
# existing date style instead of minting a new number format.

$ws.Range("A32:D43").Copy()
$ws.Paste($ws.Range("A45:D56"))

$ws.Range("F32:H43").Copy()
$ws.Paste($ws.Range("F45:H56"))

$ws.Range("E50").Value2 = $ws.Range("E37").Value2
$ws.Range("E51").Value2 = $ws.Range("E38").Value2
$ws.Range("E52").Value2 = $ws.Range("E39").Value2
$ws.Range("E53").Value2 = $ws.Range("E40").Value2
$ws.Range("E54").Value2 = $ws.Range("E41").Value2

$ws.Range("E37:E41").Copy()
$ws.Range("E50:E54").PasteSpecial(-4122)

# Make REPORT the active tab (was Data).
$ws.Activate()
